$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture values that need to move before structural edits ---
$profName    = $ws.Range("B18").Value2
$metodoVal   = $ws.Range("B19").Value2
$criterioVal = $ws.Range("B20").Value2
$normaVal    = $ws.Range("B21").Value2

# --- fix Objetivos (row 10): wrong value (was professor name) -> real objectives text ---
$objetivosTxt = @"
O objetivo da presente disciplina é introduzir os alunos no sistema normativo ambiental, conhecendo os princípios fundamentais do Direito Ambiental, sendo também capazes de analisar alguns dos instrumentos da Política Nacional de Meio Ambiente e discutir aspectos da legislação protetora dos recursos ambientais.
"@
$ws.Range("B10").Value = $objetivosTxt
$ws.Range("C10").Value = $objetivosTxt

# --- insert a new row at 13 for "Docentes responsáveis:" value (professor name) ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()
$ws.Range("B13").Value = $profName
$ws.Range("C13").Value = $profName
$ws.Range("B14:C14").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# --- fix "Programa resumido:" (now row 14): wrong value (was "Semestral") -> short syllabus pt text ---
$programaResumidoTxt = @"
Direitos ambiental constitucional; política nacional do meio ambiente
"@
$ws.Range("B14").Value = $programaResumidoTxt
$ws.Range("C14").Value = $programaResumidoTxt

# --- fix "Programa:" (now row 16): wrong value (was a date) -> full syllabus pt text ---
$programaTxt = @"
Conceitos básicos; Princípios fundamentais de direito ambiental; evolução histórica da legislação ambiental brasileira; política nacional do meio ambiente ; código florestal brasileiro; política nacional de recursos hídricos; lei dos crimes ambientais; sistema nacional de unidades de conservação; tutela administrativa, civil e processual do meio ambiente; estudos de caso com aplicação da legislação ambiental vigente e necessária para o licenciamento de empreendimento em diversos estados brasileiros.
"@
$ws.Range("B16").Value = $programaTxt
$ws.Range("C16").Value = $programaTxt

# --- fix "Método:" / "Critério:" / "Norma de recuperação:" (rows 19-21): values were each one row off ---
$ws.Range("B19").Value = $metodoVal
$ws.Range("C19").Value = $metodoVal
$ws.Range("B20").Value = $criterioVal
$ws.Range("C20").Value = $criterioVal
$ws.Range("B21").Value = $normaVal
$ws.Range("C21").Value = $normaVal

# --- fix "Bibliografia:" (now row 22): wrong value (was norma de recuperação text) -> real bibliography ---
$bibliografiaTxt = @"
ALENZA G. J. F. Manual de Derecho Ambiental. Universidad Pública de Navarra, 2001. 
ANTUNES, P. B. Dano Ambiental: uma abordagem conceitual. Rio de Janeiro, Editora Lumen Juris, 2000.
BUSTAMANTE A. J. Derecho Ambiental, Editorial Abeledo-perrot, Buenos Aires.
COSTA JR., P. J. Direito Penal Ecológico. Rio de Janeiro, Forense Universitária, 1996.
COSTA NETO, N. D. C.; BELLO FILHO, N. B.; e CASTRO E COSTA, F. D. Crimes e Infrações Administrativas Ambientais. Brasília: Brasília Jurídica, 2000. 
DICIONÁRIO DE DIREITO AMBIENTAL - Terminologia das Leis do Meio Ambiente. Maria da Graça Krieger, Anna Maria Becker Maciel, João Carlos de Carvalho Rocha, Maria José Bocorny Finatto e Cleci Regina Bevilacqua. Editora Universidade/UFRGS.
FIORILLO, C. A. P.; e RODRIGES, M. A. Manual de Direito Ambiental e legislação aplicável. São Paulo, Max Limonad, 1997.
FIORILLO, C. A. P.; e RODRIGES, M. A. Direito Ambiental e Patrimônio Genético. Belo Horizonte: Del Rey, 1996.
FREITAS, V. P. Águas - Aspectos Jurídicos e Ambientais. Curitiba, Juruá, 2000.
FREITAS, V. P. Direito Administrativo e Meio Ambiente. Curitiba, Juruá, 1993.
GOMES, C. L. S. P. Crimes Contra o Meio Ambiente: responsabilidade e sanção penal. 2ª edição, São Paulo, Editora Juarez de Oliveira, 1999.
MACHADO, P. A. L. Direito Ambiental Brasileiro. 8ª Edição, Revista, atualizada e ampliada, São Paulo, Malheiros Editores, 2001.
MORAES, A. Direito Constitucional.. 7ª ed. revista, ampliada e atualizada, com a EC n.º 24/99 - São Paulo, Atlas, 2000.
"@
$ws.Range("B22").Value = $bibliografiaTxt
$ws.Range("C22").Value = $bibliografiaTxt

Write-Output "done"
